# Atualiza os valores da coluna B (saldo_serie) nas linhas 3 a 49
# conforme "melhoria de gráficos e organização das bases"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(3, -67, -111, -108, -40, 59, 69, 96, 108, 19, -18, 54, 25, -8, -1, 0, 32, 28, 59, 76, 15, 40, -31, 73, 70, 45, 109, 24, 48, -5, 128, 99, 10, 201, -19, 60, 33, 145, 52, 139, 40, 23, -67, 137, 96, 61, -40)

$startRow = 3
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
